$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("L2").Value = 3322
$ws.Range("L3").Value = 3437
$ws.Range("I4").Value = 1839
$ws.Range("L4").Value = 863
$ws.Range("L5").Value = 198
$ws.Range("L6").Value = 3030
$ws.Range("I7").Value = 26307
$ws.Range("L7").Value = 10850

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("L2").Value = 204
$ws.Range("L3").Value = 226
$ws.Range("L6").Value = 194
$ws.Range("L7").Value = 692

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("L2").Value = 85
$ws.Range("L6").Value = 58
$ws.Range("L7").Value = 256

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("L2").Value = 144
$ws.Range("L3").Value = 153
$ws.Range("L7").Value = 502

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("L3").Value = 49
$ws.Range("L7").Value = 152

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("L2").Value = 116
$ws.Range("L3").Value = 117
$ws.Range("L7").Value = 389

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("L2").Value = 48
$ws.Range("L7").Value = 181

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("L2").Value = 93
$ws.Range("L6").Value = 86
$ws.Range("L7").Value = 365
$ws.Range("L8").Value = 692
$ws.Range("L11").Value = 184
$ws.Range("L14").Value = 54
$ws.Range("L19").Value = 305
$ws.Range("L20").Value = 276
$ws.Range("L21").Value = 32
$ws.Range("L27").Value = 103
$ws.Range("L29").Value = 587
$ws.Range("L33").Value = 502
$ws.Range("L37").Value = 389
$ws.Range("L42").Value = 346
$ws.Range("L51").Value = 136
$ws.Range("L55").Value = 104
$ws.Range("L57").Value = 40
$ws.Range("I63").Value = 263
$ws.Range("L63").Value = 35
$ws.Range("L64").Value = 72
$ws.Range("L67").Value = 391
$ws.Range("L68").Value = 36
$ws.Range("L76").Value = 154
$ws.Range("L77").Value = 65
$ws.Range("L78").Value = 137
$ws.Range("L80").Value = 33
$ws.Range("L83").Value = 256
$ws.Range("L84").Value = 108
$ws.Range("L85").Value = 550
$ws.Range("L91").Value = 154
$ws.Range("L92").Value = 33
$ws.Range("L94").Value = 130
$ws.Range("L95").Value = 152
$ws.Range("L96").Value = 106
$ws.Range("L99").Value = 181
$ws.Range("I101").Value = 26307
$ws.Range("L101").Value = 10850

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("L6").Value = 89
$ws.Range("L7").Value = 391

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("L3").Value = 43
$ws.Range("L4").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("L7").Value = 108

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("L2").Value = 175
$ws.Range("L3").Value = 225
$ws.Range("L4").Value = 29
$ws.Range("L6").Value = 150
$ws.Range("L7").Value = 587

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("L2").Value = 105
$ws.Range("L7").Value = 305

$ws = $wb.Worksheets.Item("River North")
$ws.Range("L5").Value = 6
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item("Bridgeport")
$ws.Range("L2").Value = 22
$ws.Range("L7").Value = 54

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 86

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("L3").Value = 111
$ws.Range("L7").Value = 346

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("L2").Value = 40
$ws.Range("L7").Value = 137

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 104

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("L6").Value = 26
$ws.Range("L7").Value = 106

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("L2").Value = 54
$ws.Range("L7").Value = 154

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Range("L3").Value = 7
$ws.Range("L7").Value = 32

$ws = $wb.Worksheets.Item("Near South Side")
$ws.Range("L4").Value = 8
$ws.Range("L7").Value = 72

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("L2").Value = 90
$ws.Range("L7").Value = 276

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("L2").Value = 117
$ws.Range("L7").Value = 365

$ws = $wb.Worksheets.Item("West Loop")
$ws.Range("L4").Value = 18
$ws.Range("L6").Value = 48
$ws.Range("L7").Value = 130

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("L2").Value = 70
$ws.Range("L6").Value = 44
$ws.Range("L7").Value = 184

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("L2").Value = 31
$ws.Range("L6").Value = 25
$ws.Range("L7").Value = 93

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("L6").Value = 11
$ws.Range("L7").Value = 33

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("L4").Value = 14
$ws.Range("L7").Value = 103

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("L2").Value = 41
$ws.Range("L4").Value = 19
$ws.Range("L6").Value = 32
$ws.Range("L7").Value = 136

$ws = $wb.Worksheets.Item("North Park")
$ws.Range("L2").Value = 10
$ws.Range("L3").Value = 12
$ws.Range("L7").Value = 36

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("L6").Value = 13
$ws.Range("L7").Value = 40

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("L3").Value = 224
$ws.Range("L5").Value = 11
$ws.Range("L6").Value = 117
$ws.Range("L7").Value = 550

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("L3").Value = 27
$ws.Range("L7").Value = 65

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("L6").Value = 16
$ws.Range("L7").Value = 33
